$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New row 6: duplicate of row 2 but with Number = "Demo 6" and the
#     (typo'd) Budget Level value "Approtionment" ---
$ws.Range("B6").Value2 = "Demo 6"
$ws.Range("A6").Value2 = "Approtionment"
$ws.Range("C6").Value2 = $ws.Range("C2").Value2
$ws.Range("C6").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("D6").Value2 = $ws.Range("D2").Value2
$ws.Range("E6").Value2 = $ws.Range("E2").Value2
$ws.Range("F6").Value2 = $ws.Range("F2").Value2
$ws.Range("F6").NumberFormat = $ws.Range("F2").NumberFormat
$ws.Range("G6").Value2 = $ws.Range("G2").Value2
$ws.Range("H6").Value2 = $ws.Range("H2").Value2
$ws.Range("I6").Value2 = $ws.Range("I2").Value2
$ws.Range("J6").Value2 = $ws.Range("J2").Value2
$ws.Range("K6").Value2 = $ws.Range("K2").Value2
$ws.Range("L6").Value2 = $ws.Range("L2").Value2
$ws.Range("M6").Value2 = $ws.Range("M2").Value2
$ws.Range("N6").Value2 = $ws.Range("N2").Value2
$ws.Range("O6").Value2 = $ws.Range("O2").Value2
$ws.Range("P6").Value2 = $ws.Range("P2").Value2
$ws.Range("Q6").Value2 = $ws.Range("Q2").Value2
$ws.Range("R6").Value2 = $ws.Range("R2").Value2
$ws.Range("S6").Value2 = $ws.Range("S2").Value2
$ws.Range("T6").Value2 = $ws.Range("T2").Value2
$ws.Range("U6").Value2 = $ws.Range("U2").Value2
$ws.Range("V6").Value2 = $ws.Range("V2").Value2
$ws.Range("W6").Value2 = $ws.Range("W2").Value2
$ws.Range("X6").Value2 = $ws.Range("X2").Value2

# --- New row 7: duplicate of row 3 but with Number = "Demo 6" ---
$ws.Range("A7").Value2 = $ws.Range("A3").Value2
$ws.Range("B7").Value2 = "Demo 6"
$ws.Range("C7").Value2 = $ws.Range("C3").Value2
$ws.Range("C7").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D7").Value2 = $ws.Range("D3").Value2
$ws.Range("E7").Value2 = $ws.Range("E3").Value2
$ws.Range("F7").Value2 = $ws.Range("F3").Value2
$ws.Range("F7").NumberFormat = $ws.Range("F3").NumberFormat
$ws.Range("G7").Value2 = $ws.Range("G3").Value2
$ws.Range("H7").Value2 = $ws.Range("H3").Value2
$ws.Range("I7").Value2 = $ws.Range("I3").Value2
$ws.Range("J7").Value2 = $ws.Range("J3").Value2
$ws.Range("K7").Value2 = $ws.Range("K3").Value2
$ws.Range("L7").Value2 = $ws.Range("L3").Value2
$ws.Range("M7").Value2 = $ws.Range("M3").Value2
$ws.Range("N7").Value2 = $ws.Range("N3").Value2
$ws.Range("O7").Value2 = $ws.Range("O3").Value2
$ws.Range("P7").Value2 = $ws.Range("P3").Value2
$ws.Range("Q7").Value2 = $ws.Range("Q3").Value2
$ws.Range("R7").Value2 = $ws.Range("R3").Value2
$ws.Range("S7").Value2 = $ws.Range("S3").Value2
$ws.Range("T7").Value2 = $ws.Range("T3").Value2
$ws.Range("U7").Value2 = $ws.Range("U3").Value2
$ws.Range("V7").Value2 = $ws.Range("V3").Value2
$ws.Range("W7").Value2 = $ws.Range("W3").Value2
$ws.Range("X7").Value2 = $ws.Range("X3").Value2

# Fix the typo in A2: "Apportionment" -> "Approtionment" (reuses the shared
# string already created above for A6)
$ws.Range("A2").Value2 = "Approtionment"

# Update the active selection so the file reflects the saved view state
$ws.Range("C1").Select()
